# 8/4 thermo computer changes
# Sets the Study_Batch_YYMMDD (column E) value for the rows belonging to
# the 8/4 ("250804") thermo run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$value = "CLUTEST_02_250804"

$rows = @()
for ($base = 4; $base -le 804; $base += 10) {
    $rows += $base
    $rows += ($base + 1)
}
$rows += 813

foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = $value
}
